$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two new logbook rows (row 6 and row 7) describing the work done:
# - Project organization / constraints & elements of risk / deliverables sections
# - Intro to video processing with Matlab

$ws.Range("A6").Value = "Laurent"
$ws.Range("B6").Value = "Introduction to video processing with Matlab"
$ws.Range("C6").Value = "Started to learn how to process videos with Matlab"
$ws.Range("D6").Value = 17.3
$ws.Range("E6").Value = "1h"

$ws.Range("A7").Value = "Laurent"
$ws.Range("B7").Value = "Initial plan"
$ws.Range("C7").Value = "Added: project organization, constraints & elem of risks and deliverables parts"
$ws.Range("D7").Value = 19.3
$ws.Range("E7").Value = "2h"

# Move the active selection to C21, matching where the author left off editing.
$ws.Range("C21").Select() | Out-Null
